$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2) target depth data updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 239
$wsOff.Range("C2").Value = 149
$wsOff.Range("D2").Value = 51
$wsOff.Range("E2").Value = 18
$wsOff.Range("F2").Value = 7

# DEF sheet - Home row (row 2) target depth data updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 118
$wsDef.Range("C2").Value = 83
$wsDef.Range("D2").Value = 29
$wsDef.Range("E2").Value = 14
